$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add two new test-case rows (1.6 "Test navigation to login page" and
# 1.7 "Test logout") into the previously-blank rows 8 and 9, re-using the
# same look-and-feel as the existing test-case rows above them.
# ---------------------------------------------------------------------------

# Row 8 ("#" cell looks like the other numbered first-column cells, e.g. A2)
$ws.Range("A2").Copy()
$ws.Range("A8").PasteSpecial(-4122)   # xlPasteFormats

# Row 8, columns B:G look like the body cells of the existing rows (e.g. row 7)
$ws.Range("B7:G7").Copy()
$ws.Range("B8:G8").PasteSpecial(-4122)

# Row 9 looks like row 3..7 (decimal-numbered "#" column, e.g. A3)
$ws.Range("A3:G3").Copy()
$ws.Range("A9:G9").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---- Row 8 values: Test Case 1.6 ----
$ws.Range("A8").Value = 1.6
$ws.Range("B8").Value = "Test navigation to login page"
$ws.Range("C8").Value = "To test if the login page can be navigated to successfully."
$ws.Range("D8").Value = "-"
$ws.Range("E8").Value = "The user will be brought to the login page."
$ws.Range("F8").Value = "-"
$ws.Range("G8").Value = "Fail"

# ---- Row 9 values: Test Case 1.7 ----
$ws.Range("A9").Value = 1.7
$ws.Range("B9").Value = "Test logout"
$ws.Range("C9").Value = "To check if user is able to logout."
$ws.Range("D9").Value = "-"
$ws.Range("E9").Value = "The user will be brought to the logged out page."
$ws.Range("F9").Value = "-"
$ws.Range("G9").Value = "Fail"

# ---------------------------------------------------------------------------
# Selection / view bookkeeping to mirror the saved workbook state.
# ---------------------------------------------------------------------------
$ws.Range("E7").Select()
$excel.ActiveWindow.ScrollRow = 6
